# Add a new worksheet "2024-09-21" at the end of the workbook.
# It starts life as a duplicate of the "2024-07-25" sheet (same report
# layout/header row), then the per-doctor M/F/R/Total tally columns
# (E:H) are cleared out for the new reporting period - only the
# header text in row 1 is kept, rows 2-22 keep just the ID / doctor /
# specialization / assistant columns (A:D).

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2024-07-25")

# Duplicate the template sheet; the copy is inserted right after it.
$template.Copy($null, $template)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2024-09-21"

# Clear the M / F / R / Total values for every doctor row, leaving the
# column headers in row 1 untouched.
$newSheet.Range("E2:H22").ClearContents()

# Restore the first sheet as the active one (keeps the workbook's
# selection state the same as before the edit).
$wb.Worksheets.Item(1).Activate()
